$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B2 held the "Programs" query. The website column was reworked to
# derive its value from prg.program_link / prg.program_acronym instead of
# plain prg.website, via a new CASE expression.
$newProgramsQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Pancreas Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

$cell = $ws.Range("B2")
$cell.Value = $newProgramsQuery

# Re-apply the wrap-text formatting explicitly so the edit mirrors the
# font/style bookkeeping Excel performs when a wrapped cell is retyped.
$cell.Font.Name = "Calibri"
$cell.Font.Size = 12
$cell.WrapText = $true

# Move the active selection to C3, matching where the author left off.
$ws.Range("C3").Select()
